$wb = $excel.ActiveWorkbook

# --- SE sheet: remove the stray test-data row, add a blank (quote-prefixed) row 4 ---
$wsSE = $wb.Worksheets.Item("SE")
$wsSE.Rows.Item(2).Delete()

# --- TE sheet: remove the stray test-data row entirely ---
$wsTE = $wb.Worksheets.Item("TE")
$wsTE.Rows.Item(2).Delete()

# --- BE sheet: remove the stray test-data row entirely ---
$wsBE = $wb.Worksheets.Item("BE")
$wsBE.Rows.Item(2).Delete()

# --- DT sheet: add the new column headers, widen ExamTitle column, add blank row ---
$wsDT = $wb.Worksheets.Item("DT")
$wsDT.Range("A1").Value = "Year"
$wsDT.Range("B1").Value = "Department"
$wsDT.Range("C1").Value = "ExamTitle"
$wsDT.Range("D1").Value = "NumSub"
$wsDT.Columns.Item(3).ColumnWidth = 33.6

# Blank, but formatted (quote-prefixed) row 4 under SE's header
$wsSE.Range("A4:D4").Value = "'"
$wsSE.Range("A4:D4").Value = ""

# Blank, but formatted (quote-prefixed) row 5 under DT's header
$wsDT.Range("A5:D5").Value = "'"
$wsDT.Range("A5:D5").Value = ""

# --- Update the recorded selections on each sheet ---
$wsSE.Range("A1:D5").Select()
$wsTE.Range("A2:D3").Select()
$wsBE.Range("A2:D2").Select()
$wsDT.Range("A2:D6").Select()

# Leave FE as the active sheet/tab with its own updated selection
$wsFE = $wb.Worksheets.Item("FE")
$wsFE.Range("A2").Select()
